$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Rewrite the whole "Content Placeholder 2" text so that the third bullet
# ("Completely variable in size (in 512 bytes chunks)") loses the
# " (in 512 bytes chunks)" qualifier and becomes its own short line, and
# "And an index to the log records (slot array)" becomes an independent
# paragraph again (same text, new paragraph boundary/formatting state).
$newText = "As expected, starts with a header" + [char]13 + `
           "Then a series of log records" + [char]13 + `
           "Completely variable in size" + [char]13 + `
           "And an index to the log records (slot array)"
$tr.Text = $newText

# Split "Completely variable in size" so "in " becomes its own run,
# matching "Completely variable " / "in " / "size" as three separate runs.
$full = $tr.Text
$idx = $full.IndexOf("in size")
$sub = $tr.Characters($idx + 1, 3)
$sub.Text = "in "
